$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test_value (column C) for ls_max (row 15) and ls_step (row 16)
$ws.Range("C15").Value = 360
$ws.Range("C16").Value = 120

# Move the active selection to C6 (was C12)
$ws.Range("C6").Select()
